# GetToGrid_Week1: replace the placeholder volunteer-name grid (columns D:I,
# rows 2-62) with the output of the "create grid" assignment algorithm -
# each cell becomes "<Name> <Gender>['<tag>', ...]" (or 'NONE FOUND' where
# no eligible volunteer was found). Columns A-C (index/activity/time) are
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'Nono Male['''']'
$ws.Range("F2").Value = 'Sea Female['''']'
$ws.Range("G2").Value = 'GoGo Female['''']'
$ws.Range("H2").Value = 'Surf Co Male['''']'
$ws.Range("I2").Value = 'Bow Wow Female['''']'
$ws.Range("E3").Value = 'Goobie Female['''']'
$ws.Range("F3").Value = 'Indi Male['''']'
$ws.Range("G3").Value = 'Goobie Female['''']'
$ws.Range("H3").Value = 'Indi Male['''']'
$ws.Range("I3").Value = 'Nono Male['''']'
$ws.Range("E4").Value = 'Indi Male['''']'
$ws.Range("I4").Value = 'Chicken Female['''']'
$ws.Range("E5").Value = 'Jaws Male['''']'
$ws.Range("I5").Value = 'Burning Bush Male['''']'
$ws.Range("E6").Value = 'Dad Female[''Female'', ''Leadership'']'
$ws.Range("F6").Value = 'Dad Female[''Female'', ''Leadership'']'
$ws.Range("G6").Value = 'Dad Female[''Female'', ''Leadership'']'
$ws.Range("H6").Value = 'Dad Female[''Female'', ''Leadership'']'
$ws.Range("I6").Value = 'Dad Female[''Female'', ''Leadership'']'
$ws.Range("E7").Value = 'Smiles Female[''Female'']'
$ws.Range("F7").Value = 'Goobie Female[''Female'']'
$ws.Range("G7").Value = 'T-Whisk Female[''Female'']'
$ws.Range("H7").Value = 'Unicorn Female[''Female'']'
$ws.Range("I7").Value = 'Goobie Female[''Female'']'
$ws.Range("I9").Value = 'Bow Wow Female[''Ropes'', ''Non-program'']'
$ws.Range("I10").Value = 'T-Whisk Female[''Ropes'', ''Non-program'']'
$ws.Range("I11").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("I12").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("F13").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("E14").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("F14").Value = 'Burning Bush Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("E15").Value = 'Burning Bush Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("F15").Value = 'Chicken Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("E16").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("F16").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("E17").Value = 'Captain Male[''Non-program'']'
$ws.Range("F17").Value = 'Sea Female[''Non-program'']'
$ws.Range("H17").Value = 'Laddy Male[''Non-program'']'
$ws.Range("I17").Value = 'Burning Bush Male[''Non-program'']'
$ws.Range("E18").Value = 'Bow Wow Female[''Non-program'']'
$ws.Range("F18").Value = 'Indi Male[''Non-program'']'
$ws.Range("H18").Value = 'Burning Bush Male[''Non-program'']'
$ws.Range("I18").Value = 'Sea Female[''Non-program'']'
$ws.Range("F20").Value = 'Laddy Male[''Male'', ''Non-program'']'
$ws.Range("F21").Value = 'Burning Bush Male[''Male'', ''Non-program'']'
$ws.Range("H22").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("H23").Value = 'Burning Bush Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("E24").Value = 'Tross Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("F24").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("G24").Value = 'Indi Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("H24").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("I24").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("E25").Value = 'Laddy Male[''Lifegaurd'', ''Non-program'']'
$ws.Range("F25").Value = 'Sea Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("G25").Value = 'Bow Wow Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("H25").Value = 'Sea Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("I25").Value = 'Bonez Female[''Lifegaurd'', ''Non-program'']'
$ws.Range("D26").Value = 'Nodder Male[''Lifegaurd'']'
$ws.Range("E26").Value = 'Pizza Female[''Lifegaurd'']'
$ws.Range("F26").Value = 'Chicken Female[''Lifegaurd'']'
$ws.Range("G26").Value = 'Hoops Female[''Lifegaurd'']'
$ws.Range("H26").Value = 'Bow Wow Female[''Lifegaurd'']'
$ws.Range("I26").Value = 'Nodder Male[''Lifegaurd'']'
$ws.Range("D27").Value = 'Unicorn Female[''Lifegaurd'']'
$ws.Range("E27").Value = 'Nono Male[''Lifegaurd'']'
$ws.Range("F27").Value = 'Indi Male[''Lifegaurd'']'
$ws.Range("G27").Value = 'T-Whisk Female[''Lifegaurd'']'
$ws.Range("H27").Value = 'Unicorn Female[''Lifegaurd'']'
$ws.Range("I27").Value = 'Unicorn Female[''Lifegaurd'']'
$ws.Range("D28").Value = 'Opps Female[''Lifegaurd'']'
$ws.Range("E28").Value = 'Ups Female[''Lifegaurd'']'
$ws.Range("F28").Value = 'Nono Male[''Lifegaurd'']'
$ws.Range("G28").Value = 'Nono Male[''Lifegaurd'']'
$ws.Range("H28").Value = 'Nono Male[''Lifegaurd'']'
$ws.Range("I28").Value = 'Zombie Male[''Lifegaurd'']'
$ws.Range("D29").Value = 'Bow Wow Female[''Lifegaurd'']'
$ws.Range("E29").Value = 'Nom Nom Female[''Lifegaurd'']'
$ws.Range("F29").Value = 'Captain Male[''Lifegaurd'']'
$ws.Range("G29").Value = 'Jaws Male[''Lifegaurd'']'
$ws.Range("H29").Value = 'Zombie Male[''Lifegaurd'']'
$ws.Range("I29").Value = 'G-Poppy Female[''Lifegaurd'']'
$ws.Range("E30").Value = 'Smiles Female[''Lifegaurd'']'
$ws.Range("F30").Value = 'Bonez Female[''Lifegaurd'']'
$ws.Range("G30").Value = 'Opps Female[''Lifegaurd'']'
$ws.Range("H30").Value = 'Chicken Female[''Lifegaurd'']'
$ws.Range("I30").Value = 'Goobie Female[''Lifegaurd'']'
$ws.Range("F31").Value = 'G-Poppy Female['''']'
$ws.Range("H31").Value = 'T-Whisk Female['''']'
$ws.Range("F32").Value = 'T-Whisk Female['''']'
$ws.Range("H32").Value = 'Goobie Female['''']'
$ws.Range("E33").Value = 'Goobie Female[''Female'']'
$ws.Range("F33").Value = 'Smiles Female[''Female'']'
$ws.Range("G33").Value = 'Chicken Female[''Female'']'
$ws.Range("H33").Value = 'Opps Female[''Female'']'
$ws.Range("I33").Value = 'Sea Female[''Female'']'
$ws.Range("E34").Value = 'Opps Female[''Female'']'
$ws.Range("F34").Value = 'Goobie Female[''Female'']'
$ws.Range("G34").Value = 'Bonez Female[''Female'']'
$ws.Range("H34").Value = 'Smiles Female[''Female'']'
$ws.Range("I34").Value = 'T-Whisk Female[''Female'']'
$ws.Range("E35").Value = 'Hoops Female[''Female'']'
$ws.Range("F35").Value = 'Opps Female[''Female'']'
$ws.Range("G35").Value = 'Goobie Female[''Female'']'
$ws.Range("H35").Value = 'Hoops Female[''Female'']'
$ws.Range("I35").Value = 'Opps Female[''Female'']'
$ws.Range("E36").Value = 'T-Whisk Female[''Female'']'
$ws.Range("F36").Value = 'Bow Wow Female[''Female'']'
$ws.Range("G36").Value = 'Pizza Female[''Female'']'
$ws.Range("H36").Value = 'Pizza Female[''Female'']'
$ws.Range("I36").Value = 'Bow Wow Female[''Female'']'
$ws.Range("D37").Value = 'Zombie Male['''']'
$ws.Range("E37").Value = 'Zombie Male['''']'
$ws.Range("F37").Value = 'Zombie Male['''']'
$ws.Range("G37").Value = 'Stastro Male['''']'
$ws.Range("H37").Value = 'Captain Male['''']'
$ws.Range("I37").Value = 'Nono Male['''']'
$ws.Range("E38").Value = 'Captain Male['''']'
$ws.Range("G38").Value = 'Zombie Male['''']'
$ws.Range("I38").Value = 'Laddy Male['''']'
$ws.Range("E39").Value = 'Stastro Male['''']'
$ws.Range("G39").Value = 'Blister Male['''']'
$ws.Range("I39").Value = 'GoGo Female['''']'
$ws.Range("E40").Value = 'Chicken Female['''']'
$ws.Range("G40").Value = 'Burning Bush Male['''']'
$ws.Range("I40").Value = 'Captain Male['''']'
$ws.Range("E41").Value = 'G-Poppy Female[''Ropes'']'
$ws.Range("F41").Value = 'Blister Male[''Ropes'']'
$ws.Range("G41").Value = 'G-Poppy Female[''Ropes'']'
$ws.Range("H41").Value = 'Stastro Male[''Ropes'']'
$ws.Range("I41").Value = 'Blister Male[''Ropes'']'
$ws.Range("E42").Value = 'Laddy Male[''Ropes'']'
$ws.Range("F42").Value = 'Stastro Male[''Ropes'']'
$ws.Range("G42").Value = 'Nodder Male[''Ropes'']'
$ws.Range("H42").Value = 'G-Poppy Female[''Ropes'']'
$ws.Range("I42").Value = 'Burning Bush Male[''Ropes'']'
$ws.Range("E43").Value = 'Burning Bush Male[''Ropes'']'
$ws.Range("F43").Value = 'Pizza Female[''Ropes'']'
$ws.Range("G43").Value = 'Laddy Male[''Ropes'']'
$ws.Range("H43").Value = 'Laddy Male[''Ropes'']'
$ws.Range("I43").Value = 'Stastro Male[''Ropes'']'
$ws.Range("E44").Value = 'Bow Wow Female[''Ropes'']'
$ws.Range("F44").Value = 'Nodder Male[''Ropes'']'
$ws.Range("G44").Value = 'Bow Wow Female[''Ropes'']'
$ws.Range("H44").Value = 'Burning Bush Male[''Ropes'']'
$ws.Range("I44").Value = 'Pizza Female[''Ropes'']'
$ws.Range("E45").Value = 'Nodder Male[''Ropes'']'
$ws.Range("F45").Value = 'NONE FOUND'
$ws.Range("G45").Value = 'NONE FOUND'
$ws.Range("H45").Value = 'Nodder Male[''Ropes'']'
$ws.Range("I45").Value = 'NONE FOUND'
$ws.Range("E46").Value = 'Blister Male[''Ropes'']'
$ws.Range("G46").Value = 'NONE FOUND'
$ws.Range("I46").Value = 'NONE FOUND'
$ws.Range("E47").Value = 'GoGo Female[''Female'']'
$ws.Range("G47").Value = 'Sea Female[''Female'']'
$ws.Range("I47").Value = 'Hoops Female[''Female'']'
$ws.Range("F48").Value = 'Tross Male[''Male'']'
$ws.Range("H48").Value = 'Indi Male[''Male'']'
$ws.Range("E49").Value = 'Unicorn Female[''Female'']'
$ws.Range("F49").Value = 'Hoops Female[''Female'']'
$ws.Range("H49").Value = 'Ups Female[''Female'']'
$ws.Range("I49").Value = 'Ups Female[''Female'']'
$ws.Range("F50").Value = 'Laddy Male[''Male'']'
$ws.Range("H50").Value = 'Surf Co Male[''Male'']'
$ws.Range("F51").Value = 'Burning Bush Male[''Male'']'
$ws.Range("H51").Value = 'Tross Male[''Male'']'
$ws.Range("E52").Value = 'Indi Male[''Male'']'
$ws.Range("G52").Value = 'Indi Male[''Male'']'
$ws.Range("F54").Value = 'Surf Co Male[''Lifegaurd'']'
$ws.Range("H54").Value = 'Zombie Male[''Lifegaurd'']'
$ws.Range("F55").Value = 'Blister Male[''Lifegaurd'']'
$ws.Range("H55").Value = 'Opps Female[''Lifegaurd'']'
$ws.Range("D56").Value = 'Nono Male[''Lifegaurd'']'
$ws.Range("E56").Value = 'Ups Female[''Lifegaurd'']'
$ws.Range("D57").Value = 'Hoops Female[''Lifegaurd'']'
$ws.Range("E57").Value = 'Nodder Male[''Lifegaurd'']'
$ws.Range("D58").Value = 'T-Whisk Female[''Lifegaurd'']'
$ws.Range("E58").Value = 'GoGo Female[''Lifegaurd'']'
$ws.Range("E59").Value = 'Surf Co Male['''']'
$ws.Range("F59").Value = 'Tross Male['''']'
$ws.Range("H59").Value = 'Pizza Female['''']'
$ws.Range("D60").Value = 'Nom Nom Female['''']'
$ws.Range("E60").Value = 'Opps Female['''']'
$ws.Range("F60").Value = 'Hoops Female['''']'
$ws.Range("G60").Value = 'Hoops Female['''']'
$ws.Range("H60").Value = 'Nono Male['''']'
$ws.Range("I60").Value = 'Hoops Female['''']'
$ws.Range("D61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("E61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("F61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("G61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("H61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("I61").Value = 'Acro Male[''Leadership'', ''Male'']'
$ws.Range("D62").Value = 'Zombie Male[''Male'']'
$ws.Range("E62").Value = 'Stastro Male[''Male'']'
$ws.Range("F62").Value = 'Laddy Male[''Male'']'
$ws.Range("G62").Value = 'Nono Male[''Male'']'
$ws.Range("H62").Value = 'Laddy Male[''Male'']'
$ws.Range("I62").Value = 'Nono Male[''Male'']'
